# Auto-generated Excel COM-interop script to apply cryptos.xlsx update
# (GitHub Actions data refresh: Fri Jun 28 03:52:42 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '61.815.96'
$ws.Range('E2').Value = '  +1.30%  '

# Row 3
$ws.Range('D3').Value = '3.459.31'
$ws.Range('E3').Value = '  +2.30%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  +0.01%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '583.29'
$ws.Range('E5').Value = '  +1.55%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '147.36'
$ws.Range('E6').Value = '  +7.28%  '

# Row 7
$ws.Range('D7').Value = '3.460.02'
$ws.Range('E7').Value = '  +2.37%  '

# Row 9
$ws.Range('E9').Value = '  +1.42%  '

# Row 10
$ws.Range('E10').Value = '  +0.20%  '

# Row 11
$ws.Range('E11').Value = '  +3.51%  '

# Row 12
$ws.Range('E12').Value = '  +2.83%  '

# Row 13
$ws.Range('D13').Value = '4.050.61'
$ws.Range('E13').Value = '  +2.38%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.01'
$ws.Range('E14').Value = '  +9.37%  '

# Row 15
$ws.Range('E15').Value = '  -0.91%  '

# Row 16
$ws.Range('E16').Value = '  +1.61%  '

# Row 17
$ws.Range('D17').Value = '3.455.32'
$ws.Range('E17').Value = '  +2.18%  '

# Row 18
$ws.Range('D18').Value = '61.903.68'
$ws.Range('E18').Value = '  +1.20%  '

# Row 19
$ws.Range('E19').Value = '  +8.60%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.38'
$ws.Range('E20').Value = '  +4.03%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.59'
$ws.Range('E21').Value = '  +3.19%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '389.67'
$ws.Range('E22').Value = '  +3.43%  '

# Row 23
$ws.Range('E23').Value = '  +2.92%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '73.65'
$ws.Range('E24').Value = '  +3.81%  '

# Row 25
$ws.Range('E25').Value = '  +0.27%  '

# Row 26
$ws.Range('E26').Value = '  -0.37%  '

# Row 27
$ws.Range('E27').Value = '  -2.39%  '

# Row 28
$ws.Range('D28').Value = '3.601.10'
$ws.Range('E28').Value = '  +2.25%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.183'
$ws.Range('E29').Value = '  +0.37%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.71'
$ws.Range('E30').Value = '  +3.96%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.00'
$ws.Range('E31').Value = '  +0.32%  '

# Row 32
$ws.Range('B32').Value = 'Fetch.AI'
$ws.Range('C32').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.48'
$ws.Range('E32').Value = '  -10.56%  '

# Row 33
$ws.Range('B33').Value = 'InternetComputer(DFINITY)'
$ws.Range('C33').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '8.22'
$ws.Range('E33').Value = '  +1.79%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.19'
$ws.Range('E34').Value = '  +2.19%  '

# Row 35
$ws.Range('E35').Value = '  +0.01%  '

# Row 36
$ws.Range('E36').Value = '  +3.49%  '

# Row 37
$ws.Range('D37').Value = '3.488.09'
$ws.Range('E37').Value = '  +2.45%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '7.01'
$ws.Range('E38').Value = '  +2.92%  '

# Row 39
$ws.Range('B39').Value = 'ImmutableX'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.57'
$ws.Range('E39').Value = '  +1.86%  '

# Row 40
$ws.Range('B40').Value = 'NEARProtocol'
$ws.Range('C40').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.17'
$ws.Range('E40').Value = '  +0.64%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '167.03'

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0784'
$ws.Range('E42').Value = '  +3.65%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '27.44'
$ws.Range('E43').Value = '  +7.05%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.805'
$ws.Range('E44').Value = '  +4.11%  '

# Row 45
$ws.Range('B45').Value = 'OKB'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '42.65'
$ws.Range('E45').Value = '  +2.23%  '

# Row 46
$ws.Range('B46').Value = 'Filecoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.52'
$ws.Range('E46').Value = '  +3.95%  '

# Row 48
$ws.Range('E48').Value = '  +1.76%  '

# Row 49
$ws.Range('E49').Value = '  -1.54%  '

# Row 50
$ws.Range('D50').Value = '2.572.64'
$ws.Range('E50').Value = '  +1.57%  '

# Row 51
$ws.Range('E51').Value = '  +2.56%  '
